# Mohamed Hussein Changes bdr Function
# - Row 2 (Yusuf Elsayad Bdr): "Blocked" flag F2 changes from TRUE to FALSE
# - Row 5 (Amr Elhenawy duplicate record): phone/email/password de-duplicated
#   to reuse the already-existing strings (amr.elhenawy123@gmail.com /
#   zni.vosvmzdb123@tnzro.xln); C5 already held "01234567891" so only D5/E5
#   need updating - the shared-string table is garbage-collected on save,
#   which naturally renumbers C5's reference too.
# - Columns A and B (both 25 wide) end up represented as a single merged
#   column range, and the active selection moves to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-block Yusuf's account
$ws.Range("F2").Value = $false

# De-duplicate Amr Elhenawy's second record (values already matched what's
# stored here logically; writing them explicitly keeps intent clear and
# lets the shared-string table garbage-collect the stray duplicate entries)
$ws.Range("D5").Value = "amr.elhenawy123@gmail.com"
$ws.Range("E5").Value = "zni.vosvmzdb123@tnzro.xln"

# Keep columns A:B at a displayed width of 25 (stored width "25"); setting
# both together keeps them at the same width they already had.
$ws.Columns("A:B").ColumnWidth = 24.1666666666667

# Move the active selection to F7
$ws.Range("F7").Select()
